# Split ISIC code 05T06 into ISIC 05 and ISIC 06 on each of the three
# SoBCaICbIC-* sheets. This is done by inserting a new column before the
# existing "ISIC 05T06" column (column D) on each sheet, which pushes the
# "ISIC 05T06" header (and everything to its right) one column to the
# right. Excel auto-adjusts every relative column reference inside the
# SUMIFS formulas on the row, so the only other work needed is writing
# the two new header labels ("ISIC 05" / "ISIC 06") into the column that
# used to hold "ISIC 05T06" and the newly inserted column next to it.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "SoBCaICbIC-urbanresidential",
    "SoBCaICbIC-ruralresidential",
    "SoBCaICbIC-commercial"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new column at D, shifting the old "ISIC 05T06" column (D)
    # and everything after it one column to the right.
    $ws.Columns("D:D").Insert()

    # C1 used to hold "ISIC 05T06"; it now becomes "ISIC 05".
    $ws.Range("C1").Value = "ISIC 05"
    # The newly inserted D1 becomes "ISIC 06".
    $ws.Range("D1").Value = "ISIC 06"
}
